$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update URL and Date values ---
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B2").Value = "https://molic-avc.gabriellesantosleandro.com/ValueSet/molicavc-ethnicity-valueset"
$ws1.Range("B8").Value = "2023-08-16T00:27:03-03:00"

# --- Sheet "Include from LOINC": remove the "Mixed Ethnicity" (32625-6) concept row ---
$ws2 = $wb.Worksheets.Item("Include from LOINC")
$ws2.Rows.Item(2).Delete()
